$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new title banner, merged A1:F1, centered ---
$ws.Range("A1").Value = "Time Users Took to Complete Tasks in Seconds"
$ws.Range("A1:F1").HorizontalAlignment = -4108
$ws.Range("A1:F1").Merge()

# --- Convert the per-user task timings from text ("63s", "10.28s", ...) to numeric seconds ---
$ws.Range("B3").Value = 63
$ws.Range("C3").Value = 10.28
$ws.Range("D3").Value = 23.48
$ws.Range("E3").Value = 10.53
$ws.Range("F3").Value = 38.28

$ws.Range("B4").Value = 30.48
$ws.Range("C4").Value = 7.77
$ws.Range("D4").Value = 20.76
$ws.Range("E4").Value = 34.04
$ws.Range("F4").Value = 32.61

$ws.Range("B5").Value = 88.32
$ws.Range("C5").Value = 14.55
$ws.Range("D5").Value = 31.04
$ws.Range("E5").Value = 31.8
$ws.Range("F5").Value = 26.09

$ws.Range("B6").Value = 46.67
$ws.Range("C6").Value = 27.7
$ws.Range("D6").Value = 31.39
$ws.Range("E6").Value = 30.53
# F6 stays "DNF" (user did not finish) - left untouched

$ws.Range("B7").Value = 34.36
$ws.Range("C7").Value = 8.45
$ws.Range("D7").Value = 11.76
$ws.Range("E7").Value = 20.92
$ws.Range("F7").Value = 21.66

# --- Row 8: per-task averages ---
$ws.Range("A8").Value = "Average"
$ws.Range("B8").Formula = "=AVERAGE(B3:B7)"
$ws.Range("C8:F8").Formula = "=AVERAGE(C3:C7)"

# --- Row 18: new "Tasks Matrix" banner above the existing requirements/tasks matrix, merged F18:M18 ---
$ws.Range("F18").Value = "Tasks Matrix"
$ws.Range("F18:M18").HorizontalAlignment = -4108
$ws.Range("F18:M18").Borders.LineStyle = 1
$ws.Range("F18:M18").Merge()

$ws.Range("F18").Select()
